$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row translation (German -> English)
$ws.Range("B1").Value = "salty"
$ws.Range("C1").Value = "effort"
$ws.Range("D1").Value = "takeaway"

# Column B (salty/sweet) and Column D (takeaway/cook) per-row new values
$ws.Range("B2").Value = "salty"
$ws.Range("D2").Value = "takeaway"

$ws.Range("B3").Value = "sweet"
$ws.Range("D3").Value = "cook"

$ws.Range("B4").Value = "salty"
$ws.Range("D4").Value = "takeaway"

$ws.Range("B5").Value = "salty"
$ws.Range("D5").Value = "takeaway"

$ws.Range("B6").Value = "salty"
$ws.Range("D6").Value = "cook"

$ws.Range("B7").Value = "salty"
$ws.Range("D7").Value = "cook"

$ws.Range("B8").Value = "salty"
$ws.Range("D8").Value = "cook"

$ws.Range("B9").Value = "salty"
$ws.Range("D9").Value = "cook"

$ws.Range("B10").Value = "salty"
$ws.Range("D10").Value = "cook"

$ws.Range("B11").Value = "sweet"
$ws.Range("D11").Value = "cook"

$ws.Range("B12").Value = "sweet"
$ws.Range("D12").Value = "cook"

$ws.Range("B13").Value = "sweet"
$ws.Range("D13").Value = "cook"

$ws.Range("B14").Value = "sweet"
$ws.Range("D14").Value = "cook"

$ws.Range("B15").Value = "salty"
$ws.Range("D15").Value = "cook"

$ws.Range("E6").Select()
